$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 13.0158303464755
$ws.Range("N2").Value = 1.853964204859962
$ws.Range("O2").Value = 2.02833814451736
$ws.Range("I4").Value = 19.79629629629628
$ws.Range("N4").Value = 1.98600466835246
$ws.Range("O4").Value = 2.18975222777657
$ws.Range("I9").Value = 21.79166666666666
$ws.Range("N9").Value = 2.028520339740724
$ws.Range("O9").Value = 2.242263395092639
$ws.Range("I10").Value = 21.79166666666666
$ws.Range("I12").Value = 15.36574074074072
$ws.Range("N12").Value = 1.897690627758933
$ws.Range("O12").Value = 2.081514384587572
$ws.Range("I16").Value = 12.51681286549706
$ws.Range("N16").Value = 1.844936767548521
$ws.Range("O16").Value = 2.017393709936214
$ws.Range("I17").Value = 19.79629629629628
$ws.Range("N17").Value = 1.98600466835246
$ws.Range("O17").Value = 2.18975222777657
$ws.Range("I22").Value = 21.19907407407406
$ws.Range("N22").Value = 2.015705049109126
$ws.Range("O22").Value = 2.22640738080769
$ws.Range("I23").Value = 21.19907407407406
$ws.Range("N23").Value = 2.015705049109126
$ws.Range("O23").Value = 2.22640738080769
$ws.Range("I24").Value = 21.19907407407406
$ws.Range("I25").Value = 21.19907407407406
$ws.Range("N25").Value = 2.015705049109126
$ws.Range("O25").Value = 2.22640738080769
$ws.Range("I28").Value = 21.19907407407406
$ws.Range("N28").Value = 2.015705049109126
$ws.Range("O28").Value = 2.22640738080769
$ws.Range("I37").Value = 15.74228395061728
$ws.Range("N37").Value = 1.904889690449167
$ws.Range("O37").Value = 2.090295475371289
$ws.Range("I43").Value = 15.36574074074072
$ws.Range("N43").Value = 1.897690627758933
$ws.Range("O43").Value = 2.081514384587572
$ws.Range("I44").Value = 3.38888888888889
$ws.Range("N44").Value = 1.694051767048283
$ws.Range("O44").Value = 1.836167304537999
$ws.Range("I45").Value = 2.356481481481501
$ws.Range("N45").Value = 1.678525338046114
$ws.Range("O45").Value = 1.817698795724144
$ws.Range("I54").Value = 13.0158303464755
$ws.Range("N54").Value = 1.853964204859962
$ws.Range("O54").Value = 2.02833814451736
$ws.Range("I55").Value = 18.89814814814816
$ws.Range("N55").Value = 1.967443877059447
$ws.Range("O55").Value = 2.16691042047532
$ws.Range("I56").Value = 19.30324074074072
$ws.Range("I57").Value = 19.30324074074072
$ws.Range("N57").Value = 1.975772235794973
$ws.Range("O57").Value = 2.177153507468733
$ws.Range("I60").Value = 13.75752314814816
$ws.Range("N60").Value = 1.867546171126113
$ws.Range("O60").Value = 2.044826120875009
$ws.Range("I61").Value = 13.76976495726495
$ws.Range("N61").Value = 1.867772014163364
$ws.Range("O61").Value = 2.045100507661769
$ws.Range("I64").Value = -1.819444444444444
$ws.Range("N64").Value = 1.618523362263702
$ws.Range("O64").Value = 1.746638928617865
$ws.Range("I68").Value = 21.19907407407406
$ws.Range("N68").Value = 2.015705049109126
$ws.Range("O68").Value = 2.22640738080769
